$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change: "... New users are encouraged to take the UMLS Basics Tutorial and
# to explore the [UMLS Quick Start Guide], and other training materials."
# becomes
#         "... New users are encouraged to take the UMLS Basics Tutorial and
# to explore other training materials."
#
# i.e. the "UMLS Quick Start Guide" hyperlink (and the words "the " right
# before it plus ", and " right after it) is removed, leaving a single run
# of plain text: " and to explore other training materials."
# ---------------------------------------------------------------------------

$target = $null
foreach ($h in $d.Hyperlinks) {
    if ($h.TextToDisplay -eq "UMLS Quick Start Guide") {
        $target = $h
        break
    }
}

if ($target -ne $null) {
    $linkStart = $target.Range.Start
    $linkEnd = $target.Range.End

    # Sanity-check + delete the trailing ", and " that follows the hyperlink.
    $afterRange = $d.Range($linkEnd, $linkEnd + 6)
    if ($afterRange.Text -eq ", and ") {
        $afterRange.Delete()
    }

    # Delete the hyperlink's own text (removes the hyperlink field entirely).
    $linkRange = $d.Range($linkStart, $linkEnd)
    $linkRange.Delete()

    # Sanity-check + delete the leading "the " that precedes the hyperlink.
    $beforeRange = $d.Range($linkStart - 4, $linkStart)
    if ($beforeRange.Text -eq "the ") {
        $beforeRange.Delete()
    }
}
